$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 2 (which holds Admin/admin123)
$ws.Rows("2:3").Insert()

# Populate the newly inserted rows
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "ad123"
$ws.Range("A3").Value = "ad123"
$ws.Range("B3").Value = "admin"

# Update selection to match target state
$ws.Range("B8").Select()
